# Update NATMI LR-pair (Col1a1-Itgb1) sheet with refreshed TPM-derived
# expression / specificity values, per the "update scripts wuth new tpm" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 26.81310566666667
$ws.Range("H2").Value2 = 80.439317
$ws.Range("I2").Value2 = 0.004518206005002021
$ws.Range("J2").Value2 = 0.004518206005002021
$ws.Range("M2").Value2 = 168.1098273333333
$ws.Range("N2").Value2 = 504.329482
$ws.Range("O2").Value2 = 0.2984182258032519
$ws.Range("P2").Value2 = 0.298418225803252
$ws.Range("Q2").Value2 = 4507.546563893755
$ws.Range("R2").Value2 = 40567.91907504379
$ws.Range("S2").Value2 = 0.001348315019826302
$ws.Range("T2").Value2 = 0.001348315019826302
$ws.Range("G3").Value2 = 26.81310566666667
$ws.Range("H3").Value2 = 80.439317
$ws.Range("I3").Value2 = 0.004518206005002021
$ws.Range("J3").Value2 = 0.004518206005002021
$ws.Range("O3").Value2 = 0.2893586437755394
$ws.Range("P3").Value2 = 0.2893586437755394
$ws.Range("Q3").Value2 = 4370.703421255903
$ws.Range("R3").Value2 = 39336.33079130312
$ws.Range("S3").Value2 = 0.001307381961905883
$ws.Range("T3").Value2 = 0.001307381961905883
$ws.Range("G4").Value2 = 26.81310566666667
$ws.Range("H4").Value2 = 80.439317
$ws.Range("I4").Value2 = 0.004518206005002021
$ws.Range("J4").Value2 = 0.004518206005002021
$ws.Range("M4").Value2 = 165.99353
$ws.Range("N4").Value2 = 497.98059
$ws.Range("O4").Value2 = 0.294661504941043
$ws.Range("P4").Value2 = 0.294661504941043
$ws.Range("Q4").Value2 = 4450.802059873004
$ws.Range("R4").Value2 = 40057.21853885703
$ws.Range("S4").Value2 = 0.001331341381067553
$ws.Range("T4").Value2 = 0.001331341381067553
$ws.Range("G5").Value2 = 26.81310566666667
$ws.Range("H5").Value2 = 80.439317
$ws.Range("I5").Value2 = 0.004518206005002021
$ws.Range("J5").Value2 = 0.004518206005002021
$ws.Range("M5").Value2 = 66.22673433333334
$ws.Range("N5").Value2 = 198.680203
$ws.Range("O5").Value2 = 0.1175616254801657
$ws.Range("P5").Value2 = 0.1175616254801657
$ws.Range("Q5").Value2 = 1775.744425637928
$ws.Range("R5").Value2 = 15981.69983074135
$ws.Range("S5").Value2 = 0.000531167642202283
$ws.Range("T5").Value2 = 0.000531167642202283
$ws.Range("G6").Value2 = 5771.873535333333
$ws.Range("I6").Value2 = 0.9726032482643521
$ws.Range("J6").Value2 = 0.9726032482643523
$ws.Range("M6").Value2 = 168.1098273333333
$ws.Range("N6").Value2 = 504.329482
$ws.Range("O6").Value2 = 0.2984182258032519
$ws.Range("P6").Value2 = 0.298418225803252
$ws.Range("Q6").Value2 = 970308.6634147229
$ws.Range("R6").Value2 = 8732777.970732506
$ws.Range("S6").Value2 = 0.2902425357575277
$ws.Range("T6").Value2 = 0.2902425357575278
$ws.Range("G7").Value2 = 5771.873535333333
$ws.Range("I7").Value2 = 0.9726032482643521
$ws.Range("J7").Value2 = 0.9726032482643523
$ws.Range("O7").Value2 = 0.2893586437755394
$ws.Range("P7").Value2 = 0.2893586437755394
$ws.Range("Q7").Value2 = 940851.3777387418
$ws.Range("R7").Value2 = 8467662.399648678
$ws.Range("S7").Value2 = 0.2814311568494572
$ws.Range("T7").Value2 = 0.2814311568494572
$ws.Range("G8").Value2 = 5771.873535333333
$ws.Range("I8").Value2 = 0.9726032482643521
$ws.Range("J8").Value2 = 0.9726032482643523
$ws.Range("M8").Value2 = 165.99353
$ws.Range("N8").Value2 = 497.98059
$ws.Range("O8").Value2 = 0.294661504941043
$ws.Range("P8").Value2 = 0.294661504941043
$ws.Range("Q8").Value2 = 958093.6628435596
$ws.Range("R8").Value2 = 8622842.965592038
$ws.Range("S8").Value2 = 0.2865887368441208
$ws.Range("T8").Value2 = 0.2865887368441209
$ws.Range("G9").Value2 = 5771.873535333333
$ws.Range("I9").Value2 = 0.9726032482643521
$ws.Range("J9").Value2 = 0.9726032482643523
$ws.Range("M9").Value2 = 66.22673433333334
$ws.Range("N9").Value2 = 198.680203
$ws.Range("O9").Value2 = 0.1175616254801657
$ws.Range("P9").Value2 = 0.1175616254801657
$ws.Range("Q9").Value2 = 382252.3352301181
$ws.Range("R9").Value2 = 3440271.017071063
$ws.Range("S9").Value2 = 0.1143408188132463
$ws.Range("T9").Value2 = 0.1143408188132464
$ws.Range("G10").Value2 = 132.4457753333333
$ws.Range("H10").Value2 = 397.337326
$ws.Range("I10").Value2 = 0.02231808970163987
$ws.Range("J10").Value2 = 0.02231808970163988
$ws.Range("M10").Value2 = 168.1098273333333
$ws.Range("N10").Value2 = 504.329482
$ws.Range("O10").Value2 = 0.2984182258032519
$ws.Range("P10").Value2 = 0.298418225803252
$ws.Range("Q10").Value2 = 22265.43642231612
$ws.Range("R10").Value2 = 200388.9278008451
$ws.Range("S10").Value2 = 0.006660124732081199
$ws.Range("T10").Value2 = 0.006660124732081202
$ws.Range("G11").Value2 = 132.4457753333333
$ws.Range("H11").Value2 = 397.337326
$ws.Range("I11").Value2 = 0.02231808970163987
$ws.Range("J11").Value2 = 0.02231808970163988
$ws.Range("O11").Value2 = 0.2893586437755394
$ws.Range("P11").Value2 = 0.2893586437755394
$ws.Range("Q11").Value2 = 21589.48726703972
$ws.Range("R11").Value2 = 194305.3854033575
$ws.Range("S11").Value2 = 0.006457932167727347
$ws.Range("T11").Value2 = 0.006457932167727349
$ws.Range("G12").Value2 = 132.4457753333333
$ws.Range("H12").Value2 = 397.337326
$ws.Range("I12").Value2 = 0.02231808970163987
$ws.Range("J12").Value2 = 0.02231808970163988
$ws.Range("M12").Value2 = 165.99353
$ws.Range("N12").Value2 = 497.98059
$ws.Range("O12").Value2 = 0.294661504941043
$ws.Range("P12").Value2 = 0.294661504941043
$ws.Range("Q12").Value2 = 21985.14178116692
$ws.Range("R12").Value2 = 197866.2760305023
$ws.Range("S12").Value2 = 0.006576281898894398
$ws.Range("T12").Value2 = 0.0065762818988944
$ws.Range("G13").Value2 = 132.4457753333333
$ws.Range("H13").Value2 = 397.337326
$ws.Range("I13").Value2 = 0.02231808970163987
$ws.Range("J13").Value2 = 0.02231808970163988
$ws.Range("M13").Value2 = 66.22673433333334
$ws.Range("N13").Value2 = 198.680203
$ws.Range("O13").Value2 = 0.1175616254801657
$ws.Range("P13").Value2 = 0.1175616254801657
$ws.Range("Q13").Value2 = 8771.45117657302
$ws.Range("R13").Value2 = 78943.06058915719
$ws.Range("S13").Value2 = 0.002623750902936929
$ws.Range("T13").Value2 = 0.002623750902936929
$ws.Range("G14").Value2 = 3.326003
$ws.Range("H14").Value2 = 9.978009
$ws.Range("I14").Value2 = 0.0005604560290058679
$ws.Range("J14").Value2 = 0.000560456029005868
$ws.Range("M14").Value2 = 168.1098273333333
$ws.Range("N14").Value2 = 504.329482
$ws.Range("O14").Value2 = 0.2984182258032519
$ws.Range("P14").Value2 = 0.298418225803252
$ws.Range("Q14").Value2 = 559.1337900401487
$ws.Range("R14").Value2 = 5032.204110361338
$ws.Range("S14").Value2 = 0.000167250293816667
$ws.Range("T14").Value2 = 0.0001672502938166671
$ws.Range("G15").Value2 = 3.326003
$ws.Range("H15").Value2 = 9.978009
$ws.Range("I15").Value2 = 0.0005604560290058679
$ws.Range("J15").Value2 = 0.000560456029005868
$ws.Range("O15").Value2 = 0.2893586437755394
$ws.Range("P15").Value2 = 0.2893586437755394
$ws.Range("Q15").Value2 = 542.1592288460404
$ws.Range("R15").Value2 = 4879.433059614363
$ws.Range("S15").Value2 = 0.0001621727964489623
$ws.Range("T15").Value2 = 0.0001621727964489623
$ws.Range("G16").Value2 = 3.326003
$ws.Range("H16").Value2 = 9.978009
$ws.Range("I16").Value2 = 0.0005604560290058679
$ws.Range("J16").Value2 = 0.000560456029005868
$ws.Range("M16").Value2 = 165.99353
$ws.Range("N16").Value2 = 497.98059
$ws.Range("O16").Value2 = 0.294661504941043
$ws.Range("P16").Value2 = 0.294661504941043
$ws.Range("Q16").Value2 = 552.09497876059
$ws.Range("R16").Value2 = 4968.85480884531
$ws.Range("S16").Value2 = 0.0001651448169601499
$ws.Range("T16").Value2 = 0.0001651448169601499
$ws.Range("G17").Value2 = 3.326003
$ws.Range("H17").Value2 = 9.978009
$ws.Range("I17").Value2 = 0.0005604560290058679
$ws.Range("J17").Value2 = 0.000560456029005868
$ws.Range("M17").Value2 = 66.22673433333334
$ws.Range("N17").Value2 = 198.680203
$ws.Range("O17").Value2 = 0.1175616254801657
$ws.Range("P17").Value2 = 0.1175616254801657
$ws.Range("Q17").Value2 = 220.2703170728697
$ws.Range("R17").Value2 = 1982.432853655827
$ws.Range("S17").Value2 = 0.0000658881217800887
$ws.Range("T17").Value2 = 0.00006588812178008872
